$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.552.76"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.365.25"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.20%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "556.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.355.22"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("E9").Value = "  +0.06%  "
$ws.Range("E10").Value = "  +0.84%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.162"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "54.54"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000273"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.908.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.63%  "
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.380.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.87"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "64.483.59"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.51%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.985"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.44%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "456.72"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +12.08%  "
$ws.Range("E23").Value = "  +9.21%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.09"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.52"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.37%  "
$ws.Range("E27").Value = "  +1.31%  "
$ws.Range("E28").Value = "  +1.77%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.81"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "29.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.34%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "583.28"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.28%  "
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.70"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.82%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("E37").Value = "  -7.53%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.47"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0₃0758"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.42%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.375"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.31%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.108.44"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.80"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.52"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.70%  "
$ws.Range("E46").Value = "  +0.34%  "
$ws.Range("E47").Value = "  +0.14%  "
$ws.Range("E48").Value = "  -0.25%  "
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.33"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "135.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
